# Scout Html UI Meeting.docx - apply commit's edits
$d = $word.ActiveDocument

# 1) Update the SAVEDATE field result text "28.05.2014 11:39" -> "28.05.2014 12:52"
$d.Content.Find.Execute("28.05.2014 11:39", $true, $false, $false, $false, $false,
                         $true, 1, $false, "28.05.2014 12:52", 2)

# 2) Remove the _GoBack bookmark from the "Done" table cell paragraph
#    (it gets re-created later at the new content's last edit point)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Append a new meeting section at the end of the document, just before
#    the trailing empty paragraph / section break.
$last = $d.Paragraphs.Last
$target = $last.Range
$target.Collapse(1)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
            </w:pPr>
            <w:r>
              <w:t>04.06.2014</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="berschrift3Zwischentitel"/>
            </w:pPr>
            <w:r>
              <w:t>Teilnehmer: awe, bsh, cgu</w:t>
            </w:r>
            <w:r>
              <w:br/>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Aufzhlung"/>
            </w:pPr>
            <w:r>
              <w:t>Wir machen z.Z. noch keinen Scout-Fork f&#252;r Html UI und versuchen so lange es geht mit dem &#8222;extension&#8220;-PlugIn weiterzumachen.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Aufzhlung"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="0"/>
              </w:numPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)
